# Add Trade #35 (closed row 36) to the "All Trades" and "base_strategy" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Cells.Item(36, 1).Value = 35
    $ws.Cells.Item(36, 2).Value = "2026-02-16"
    $ws.Cells.Item(36, 3).Value = "22:55:36"
    $ws.Cells.Item(36, 4).Value = "base_strategy"
    $ws.Cells.Item(36, 5).Value = "DOWN"
    $ws.Cells.Item(36, 6).Value = 49.999998
    $ws.Cells.Item(36, 7).Value = ""
    $ws.Cells.Item(36, 8).Value = "OPEN"
    $ws.Cells.Item(36, 9).Value = 0
    $ws.Cells.Item(36, 10).Value = 0
    $ws.Cells.Item(36, 11).Value = 100
    $ws.Cells.Item(36, 12).Value = 0
    $ws.Cells.Item(36, 13).Value = 0
    $ws.Cells.Item(36, 14).Value = 0.6
    $ws.Cells.Item(36, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(36, 16).Value = ""
    $ws.Cells.Item(36, 17).Value = 0
}
